# computational speed up 3X
# auto adjust mzWindow and stepsize for gaussian peak convolution, speed up by 3X
#
# This script:
#  1. Updates the isotope-correction results on "Corrected" and "Normalized"
#     with the values produced by the re-tuned (auto mzWindow/stepsize) solver.
#  2. Adds a new "logs" worksheet (at the end of the workbook) recording the
#     solver parameters/diagnostics for the run (purity, resolution, ppm,
#     runtime).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Corrected" sheet - updated absolute intensities (C2:H6)
# ---------------------------------------------------------------------------
$corrected = $wb.Worksheets.Item("Corrected")

$corrected.Range("C2").Value = 27124615.93765815
$corrected.Range("D2").Value = 39143096.974471189
$corrected.Range("E2").Value = 26952896.289438937
$corrected.Range("F2").Value = 11356733.143402042
$corrected.Range("G2").Value = 17521967.037179444
$corrected.Range("H2").Value = 17440799.969836958

$corrected.Range("C3").Value = 4615669.5820185514
$corrected.Range("D3").Value = 5526537.2713224366
$corrected.Range("E3").Value = 4225642.3788505234
$corrected.Range("F3").Value = 8869941.0654416699
$corrected.Range("G3").Value = 18405549.062648293
$corrected.Range("H3").Value = 14791524.739797054

$corrected.Range("C4").Value = 2626110.2295785537
$corrected.Range("D4").Value = 2843239.6165311886
$corrected.Range("E4").Value = 2248276.6633403036
$corrected.Range("F4").Value = 7929278.0899684448
$corrected.Range("G4").Value = 18486321.550155826
$corrected.Range("H4").Value = 13828434.712488767

$corrected.Range("C5").Value = 1158054.1991937698
$corrected.Range("D5").Value = 1773945.3968971909
$corrected.Range("E5").Value = 1195308.3693457516
$corrected.Range("F5").Value = 5331860.8828744236
$corrected.Range("G5").Value = 13949468.279744884
$corrected.Range("H5").Value = 9937551.3175618947

$corrected.Range("C6").Value = 106100.19669463292
$corrected.Range("D6").Value = 135186.7945610942
$corrected.Range("E6").Value = 99743.593150817411
$corrected.Range("F6").Value = 1680203.6635350587
$corrected.Range("G6").Value = 5735109.2012338359
$corrected.Range("H6").Value = 2972218.1773711625

# ---------------------------------------------------------------------------
# 2. "Normalized" sheet - updated fractional abundances (C2:H6)
# ---------------------------------------------------------------------------
$normalized = $wb.Worksheets.Item("Normalized")

$normalized.Range("C2").Value = 0.76127
$normalized.Range("D2").Value = 0.79200999999999999
$normalized.Range("E2").Value = 0.77625
$normalized.Range("F2").Value = 0.32291999999999998
$normalized.Range("G2").Value = 0.23646
$normalized.Range("H2").Value = 0.29575000000000001

$normalized.Range("C3").Value = 0.12953999999999999
$normalized.Range("D3").Value = 0.11182
$normalized.Range("E3").Value = 0.12169000000000001
$normalized.Range("F3").Value = 0.25220999999999999
$normalized.Range("G3").Value = 0.24839
$normalized.Range("H3").Value = 0.25081999999999999

$normalized.Range("C4").Value = 0.073700000000000002
$normalized.Range("D4").Value = 0.057520000000000002
$normalized.Range("E4").Value = 0.064750000000000002
$normalized.Range("F4").Value = 0.22545999999999999
$normalized.Range("G4").Value = 0.24948000000000001
$normalized.Range("H4").Value = 0.23449

$normalized.Range("C5").Value = 0.032500000000000001
$normalized.Range("D5").Value = 0.035889999999999998
$normalized.Range("E5").Value = 0.034419999999999999
$normalized.Range("F5").Value = 0.15160999999999999
$normalized.Range("G5").Value = 0.18825
$normalized.Range("H5").Value = 0.16850999999999999

$normalized.Range("C6").Value = 0.00297
$normalized.Range("D6").Value = 0.0027299999999999998
$normalized.Range("E6").Value = 0.0028700000000000002
$normalized.Range("F6").Value = 0.04777
$normalized.Range("G6").Value = 0.07739
$normalized.Range("H6").Value = 0.0504

# ---------------------------------------------------------------------------
# 3. New "logs" worksheet - appended as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$logs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$logs.Name = "logs"

$logs.Range("A1").Value = "solver"
$logs.Range("B1").Value = "optcorr"

$logs.Range("A2").Value = "purity"
$logs.Range("B2").Value = 0.98999999999999999

$logs.Range("A3").Value = "resolution"
$logs.Range("B3").Value = 140000

$logs.Range("A4").Value = "ppm"
$logs.Range("B4").Value = 5

$logs.Range("A5").Value = "runtime"
$logs.Range("B5").Value = 17.6821725
